# Update the "Valid Species" sheet's alphabetical species list:
#   1. Add "Pogoniulus atroflavus" (missing species) in alphabetical order,
#      i.e. immediately before "Pogoniulus bilineatus".
#   2. Remove "Turdus philomelos" and its four subspecies entries
#      ("Turdus philomelos clarkei", "Turdus philomelos hebridensis",
#      "Turdus philomelos nataliae", "Turdus philomelos philomelos").
#
# (per commit message: "removed dates from metadata and added to
# datasets" / taxonomy clean-up while working on mptt for taxa)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid Species")

# --- 1. Insert "Pogoniulus atroflavus" right before "Pogoniulus bilineatus" ---
$ws.Rows.Item(3160).Insert()
$ws.Range("A3160").Value = "Pogoniulus atroflavus"

# --- 2. Delete the 5 "Turdus philomelos" rows (shifted down by 1 row from the insert above) ---
$firstRow = $ws.Rows.Item(4137)
$lastRow = $ws.Rows.Item(4141)
$rng = $ws.Range($firstRow, $lastRow)
$rng.EntireRow.Delete()
